$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 332 (this shifts existing row 332.. down to 333..)
$ws.Rows.Item(332).Insert()

# Populate the newly inserted row 332 with the new data record
$ws.Cells.Item(332, 1).Value = 4
$ws.Cells.Item(332, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(332, 3).Value = "Los Lagos"
$ws.Cells.Item(332, 4).Value = 44992
$ws.Cells.Item(332, 5).Value = 10
$ws.Cells.Item(332, 6).Value = 100112027
$ws.Cells.Item(332, 7).Value = "Melón"
$ws.Cells.Item(332, 8).Value = "Tuna"
$ws.Cells.Item(332, 9).Value = "Extra"
$ws.Cells.Item(332, 10).Value = 160
$ws.Cells.Item(332, 11).Value = 12000
$ws.Cells.Item(332, 12).Value = 12000
$ws.Cells.Item(332, 13).Value = 12000
$ws.Cells.Item(332, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(332, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(332, 16).Value = 1000
$ws.Cells.Item(332, 17).Value = 12
$ws.Cells.Item(332, 18).Value = "Hortaliza"
